# Staff meeting fourth version
#
# Three shape-level tweaks on the "DevOps POC" slides:
#   1. Slide 5 "TextBox 24"        - nudge position & reduce rotation.
#   2. Slide 5 "Left-Right Arrow 27" - nudge position & reduce rotation.
#   3. Slide 8 "Oval 6"            - change fill from the yellow accent
#                                     swatch (FBE700) to the background-1
#                                     theme color (white).
#
# NOTE on the literal Left/Top/Rotation numbers below: the host's Single
# (float32) precision for EMU offsets truncates (rather than rounds) when
# converting points -> EMU, so the literals are chosen so that after the
# float32 round-trip they truncate to the exact target EMU values from the
# authoritative OOXML (off x/y in EMU, 914400 EMU/in, 12700 EMU/pt).

$p = $ppt.ActivePresentation

# ---- Slide 5 ----
$s5 = $p.Slides.Item(5)

$textBox24 = $s5.Shapes.Item(19)   # "TextBox 24"
$textBox24.Rotation = 34.43005
$textBox24.Left = 323.75547792086616
$textBox24.Top = 185.92114256220472

$arrow27 = $s5.Shapes.Item(23)     # "Left-Right Arrow 27"
$arrow27.Rotation = 34.6492
$arrow27.Left = 295.27169799330704
$arrow27.Top = 201.3536605972441

# ---- Slide 8 ----
$s8 = $p.Slides.Item(8)

$oval6 = $s8.Shapes.Item(6)        # "Oval 6"
$oval6.Fill.ForeColor.SchemeColor = "bg1"
